$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update G/H values for existing rows 3-23 ---
$updates = @(
    @{ Row = 3;  G = 3216; H = 3052 },
    @{ Row = 4;  G = 3612; H = 3093 },
    @{ Row = 5;  G = 4458; H = 3114 },
    @{ Row = 6;  G = 5084; H = 3114 },
    @{ Row = 7;  G = 4977; H = 2908 },
    @{ Row = 8;  G = 4299; H = 3009 },
    @{ Row = 9;  G = 3905; H = 2925 },
    @{ Row = 10; G = 3378; H = 2978 },
    @{ Row = 11; G = 2981; H = 2932 },
    @{ Row = 12; G = 2772; H = 3049 },
    @{ Row = 13; G = 2768; H = 2810 },
    @{ Row = 14; G = 2724; H = 2822 },
    @{ Row = 15; G = 2681; H = 2860 },
    @{ Row = 16; G = 2691; H = 2806 },
    @{ Row = 17; G = 2690; H = 2891 },
    @{ Row = 18; G = 2659; H = 3063 },
    @{ Row = 19; G = 2633; H = 2835 },
    @{ Row = 20; G = 2610; H = 2855 },
    @{ Row = 21; G = 2519; H = 2848 },
    @{ Row = 22; G = 2666; H = 2835 },
    @{ Row = 23; G = 2640; H = 2944 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
}

# --- New rows 24 and 25 ---
$ws.Cells.Item(24, 6).Value = 32
$ws.Cells.Item(24, 7).Value = 2623
$ws.Cells.Item(24, 8).Value = 3103

$ws.Cells.Item(25, 6).Value = 33
$ws.Cells.Item(25, 7).Value = 3172
$ws.Cells.Item(25, 8).Value = 3114

# --- Column I: G-H difference formula for rows 3-25 ---
# Fill as one range so Excel registers the shared-formula group used by rows
# 4:25, then row 3 and row 24 are (re)written individually to match the
# original layout (row 3 stand-alone, row 24 stand-alone, row 25 rejoining
# the shared group).
$ws.Range("I3").Formula = "=G3-H3"
$ws.Range("I4:I25").Formula = "=G4-H4"
$ws.Range("I24").Formula = "=G24-H24"

# --- Update totals row 28 to extend range through row 25 ---
$ws.Range("G28").Formula = "=SUM(G3:G25)"
$ws.Range("H28").Formula = "=SUM(H3:H25)"
$ws.Range("I28").Formula = "=SUM(I3:I25)"

# --- Update view selection (no more topLeftCell="A2" freeze, select I6 instead of I24) ---
$ws.Range("I6").Select()
